$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update H2 value from 2.27 to 3.17
$ws.Range("H2").Value = 3.17

# Insert a new row 5 by copying row 4 (carries over cell styles/formatting)
$ws.Rows.Item(4).Copy()
$ws.Rows.Item(5).Insert()

# Set the new row's data: D5 = 4 (Pull #), E5 = "2C#2" (Cable Size)
$ws.Range("D5").Value = 4
$ws.Range("E5").Value = "2C#2"

# Extend the merged cells in columns A, B, C, F, G, H, I, J from row 2:4 to row 2:5
$ws.Range("A2:A4").UnMerge()
$ws.Range("B2:B4").UnMerge()
$ws.Range("C2:C4").UnMerge()
$ws.Range("F2:F4").UnMerge()
$ws.Range("G2:G4").UnMerge()
$ws.Range("H2:H4").UnMerge()
$ws.Range("I2:I4").UnMerge()
$ws.Range("J2:J4").UnMerge()

$ws.Range("A2:A5").Merge()
$ws.Range("B2:B5").Merge()
$ws.Range("C2:C5").Merge()
$ws.Range("F2:F5").Merge()
$ws.Range("G2:G5").Merge()
$ws.Range("H2:H5").Merge()
$ws.Range("I2:I5").Merge()
$ws.Range("J2:J5").Merge()
